# Add new quotes to the "group" worksheet of the quotes workbook.
# Mirrors the author manually typing 9 new quote/author rows (66-74)
# right after the existing data (which ends at row 65).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("group")

$newQuotes = @(
    @("""Food is not just fuel, it's information. It talks to your DNA and tells it what to do.""", "Dr. Mark Hyman"),
    @("""The discovery of a new dish does more for human happiness than the discovery of a new star.""", "Anthelme Brillat-Savarin"),
    @("""Food is memories.""", "José Andrés"),
    @("""Good bread is the most fundamentally satisfying of all foods; and good bread with fresh butter, the greatest of feasts.""", "James Beard"),
    @("""Food is the ingredient that binds us together.""", "Unknown"),
    @("""If you can eat with mates or friends or family, I mean, it's such a brilliant thing isn't it? If you feel really rubbish and you have a nice bit of food it makes you feel good, you know?""", "Jamie Oliver"),
    @("""Food is the most primitive form of comfort.""", "Sheila Graham"),
    @("""Food is the ultimate equalizer. It doesn't matter who you are or where you come from, everyone has to eat.""", "Unknown"),
    @("""Food is love made visible.""", "Unknown")
)

$startRow = 66

# Fill the "author" column first for every new row, then the "quote"
# column, matching how the rows were originally authored (column by
# column) so new shared-string entries line up the same way.
for ($i = 0; $i -lt $newQuotes.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $newQuotes[$i][1]
}
for ($i = 0; $i -lt $newQuotes.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newQuotes[$i][0]
}

$lastRow = $startRow + $newQuotes.Count - 1

# Update the view so the newly typed block is visible/selected, like Excel
# records after the user scrolls down and selects the freshly entered rows.
$ws.Activate()
$ws.Range("A" + $startRow + ":B" + $lastRow).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
